$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Q1 (TC-001) expected result - row 2, column E
$ws.Range("E2").Value = "1. Програма запустилась успішно 2. Час запуску програми не перевищує 2 с; використання ресурсів процесора не перевищує 20%. 3. Програма закрилась успішно 4. Час закриття програми не перевищує 2 с; використання ресурсів процесора не перевищує 20%."

# Q2 (TC-002) expected result - row 3, column E
$ws.Range("E3").Value = "1. Програма відкрилась без помилок 2. Відкрився обраний пункт меню 3. Час реакції на вибір пункта меню не перевищує 2 с."

# Q3 (TC-003) expected result - row 4, column E
$ws.Range("E4").Value = "1. Програма запустилась без помилок 2. Клітинка на яку натиснули відкрилась 3. Час реакції на натискання кнопок миші не перевищує 0.2 с. 4. Клітинка на яку натиснули позначилась прапором. 5.Час реакції на натискання кнопок миші не перевищує 0.2 с. 4."

# Last table row (row 5) was missing its own test-case id and incorrectly reused
# TC-003's id. Give it its own id: TC-004.
$ws.Range("A5").Value = "TC-004"

# Row 5 expected result text, updated with expected results per step.
$ws.Range("E5").Value = "1. Програма запустилась успішно 2. Почалась нова гра, ігрове поле оновилось 3. Час між натисканням кнопки старту нової гри та оновленням ігрового поля не перевищує 1 с."
